$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Sheet1 ("Financial Statement") has two small "Wealth Class" breakdown lists
# (rows 12-14 and rows 17-19 originally). Each list gets two additional rows
# inserted (LLA, HH, ZUT replace the old "APER"/"MENDA" pairing with a longer
# AA/MEN/LLA/HH/Fixed Asset chain), and the two existing value rows change
# their text from APER/MENDA to AA/MEN.
# ---------------------------------------------------------------------------

# --- First list: rows 12-14 (currently APER / MENDA / Fixed Asset) ---------

# Insert 3 blank rows right before row 14 ("Fixed Asset"). The new rows
# inherit row 13's format (thin box border, style index used by row13).
$ws.Rows.Item(14).Insert()
$ws.Rows.Item(14).Insert()
$ws.Rows.Item(14).Insert()

# Row 14 already has row13's border format from the insert - just set text.
$ws.Range("A14:N14").Value = ""
$ws.Range("A14").Value = "LLA"

# Row 15 needs row12's format (box without top edge) - copy/paste row 12.
$ws.Range("A12:N12").Copy()
$ws.Range("A15").Select()
$ws.Paste()
$ws.Range("A15:N15").Value = ""
$ws.Range("A15").Value = "HH"

# Row 16 needs the "Fixed Asset" row's format - that row is now at position
# 17 (after the 3 inserts). Copy/paste it, then replace the text.
$ws.Range("A17:N17").Copy()
$ws.Range("A16").Select()
$ws.Paste()
$ws.Range("A16:N16").Value = ""
$ws.Range("A16").Value = "ZUT"

# Update the two pre-existing rows' labels (APER -> AA, MENDA -> MEN).
$ws.Range("A12").Value = "AA"
$ws.Range("A13").Value = "MEN"

$excel.CutCopyMode = 0

# --- Second list: now at rows 20 ("APER"), 21 ("MENDA"), 22 ("Fixed Asset")
#     (shifted down by 3 from the original 17/18/19 because of the inserts
#     above). -------------------------------------------------------------

# Insert 3 blank rows right before row 22 ("Fixed Asset"). The new rows
# inherit row 21's format (same border style as row 18 originally had).
$ws.Rows.Item(22).Insert()
$ws.Rows.Item(22).Insert()
$ws.Rows.Item(22).Insert()

# Row 22: column A keeps the inherited style from row 21; columns B:N get
# re-formatted to match row 12's body style (applyBorder, no top edge).
$ws.Range("A22").Value = "LLA"
$ws.Range("B12:N12").Copy()
$ws.Range("B22").Select()
$ws.Paste()
$ws.Range("B22:N22").Value = ""

# Row 23 needs the thin full-box style (same as row 13 / row 14).
$ws.Range("A14:N14").Copy()
$ws.Range("A23").Select()
$ws.Paste()
$ws.Range("A23:N23").Value = ""
$ws.Range("A23").Value = "HH"

# Row 24 needs the "Fixed Asset" row style - donor is row 16 built above.
$ws.Range("A16:N16").Copy()
$ws.Range("A24").Select()
$ws.Paste()
$ws.Range("A24:N24").Value = ""
$ws.Range("A24").Value = "ZUT"

# Update the two pre-existing rows' labels (APER -> AA, MENDA -> MEN).
$ws.Range("A20").Value = "AA"
$ws.Range("A21").Value = "MEN"

# Final "Fixed Asset" row (now at 25) switches from its old style to the
# "Fixed Asset" box style used elsewhere (copy from row 16/24's style).
$ws.Range("A16:N16").Copy()
$ws.Range("A25").Select()
$ws.Paste()
$ws.Range("A25:N25").Value = ""
$ws.Range("A25").Value = "Fixed Asset"

$excel.CutCopyMode = 0

$wb.Save()
